# Implement hierarchical TreeView and editing for Frais divers
# - "Entretien" sheet: row 2 description edited and amount column repurposed
#   to hold a (text) date value instead of a numeric amount.
# - "frais divers" sheet: row 2 amount updated from 1000 to 1500.

$wb = $excel.ActiveWorkbook

# --- Sheet "Entretien" (row 2: A2 description, B2 now a text date) ---
$wsEntretien = $wb.Worksheets.Item("Entretien")
$wsEntretien.Range("A2").Value = "Modified Description"

# Force B2 to stay a literal text string ("2024-01-15") instead of being
# auto-converted to a date serial number, then reset the style so no
# extra number-format styling is attached to the cell.
$wsEntretien.Range("B2").Value = "'2024-01-15"
$wsEntretien.Range("B2").Style = "Normal"

# --- Sheet "frais divers" (row 2: B2 amount 1000 -> 1500) ---
$wsFraisDivers = $wb.Worksheets.Item("frais divers")
$wsFraisDivers.Range("B2").Value = 1500
